# Add a new forecast-origin column (AH) with header date "2020-05-05"
# and a new observed-date row (46) for "2020-05-19" to both the
# "cases" and "deaths" worksheets, plus a handful of updated/filled-in
# diagonal forecast values.

$wb = $excel.ActiveWorkbook

# New diagonal values for column AH (rows 33-46), keyed by sheet name.
$casesAH = @{
    33 = 123085; 34 = 130735; 35 = 137934; 36 = 143674; 37 = 149970;
    38 = 156123; 39 = 162697; 40 = 167780; 41 = 172715; 42 = 176630;
    43 = 180773; 44 = 184671; 45 = 188755; 46 = 192203
}
$deathsAH = @{
    33 = 8364; 34 = 8829; 35 = 9268; 36 = 9660; 37 = 10015;
    38 = 10366; 39 = 10710; 40 = 10996; 41 = 11281; 42 = 11535;
    43 = 11769; 44 = 11984; 45 = 12201; 46 = 12421
}

$sheetConfigs = @(
    @{ Name = "cases";  B31New = 107780; B32New = 114715; AH = $casesAH },
    @{ Name = "deaths"; B31New = 7321;   B32New = 7921;   AH = $deathsAH }
)

foreach ($cfg in $sheetConfigs) {
    $ws = $wb.Worksheets.Item($cfg.Name)

    # --- New column AH -------------------------------------------------
    # Header cell AH1 carries the forecast-origin date "2020-05-05",
    # the same text already used by A32 / AG-column dates further down.
    $ws.Range("AH1").NumberFormat = "@"
    $ws.Range("AH1").Value = "2020-05-05"
    $ws.Range("AH1").Style = "Normal"

    # Rows 2-32: AH is present but empty for all of them.
    $ws.Range("AH2:AH32").Style = "Normal"

    # Rows 33-46: AH carries the new diagonal forecast numbers.
    foreach ($r in 33..46) {
        $ws.Cells.Item($r, 34).Value = $cfg.AH[$r]
    }

    # --- Existing column B updates --------------------------------------
    $ws.Range("B31").Value = $cfg.B31New
    $ws.Range("B32").Value = $cfg.B32New

    # --- New row 46 (observed date 2020-05-19) --------------------------
    # A46 gets the new date label; B46..AG46 stay empty (only AH46, set
    # above, carries a value).
    $ws.Range("B46:AG46").Style = "Normal"
    $ws.Range("A46").NumberFormat = "@"
    $ws.Range("A46").Value = "2020-05-19"
    $ws.Range("A46").Style = "Normal"
}
